$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated vm_pu results for the case with 380 kV (row index = data row, column letters per sheet)
$data = @{
    2 = @{ "B"=1.02; "C"=1.047755892524305; "D"=1.0521603800662; "E"=1.060573447651939; "F"=1.066268244441687; "I"=1.040323915334692; "J"=1.052803269752584; "K"=1.05490988118869; "L"=1.0632998589796; "M"=1.068979250086863; "N"=1.054298371153855 }
    3 = @{ "B"=1.02; "C"=1.048786105183069; "D"=1.052933811906198; "E"=1.06154080503325; "F"=1.067204813968297; "I"=1.040504152781486; "J"=1.053481687299104; "K"=1.055496204448358; "L"=1.064081282271395; "M"=1.069731077108201; "N"=1.054977752131107 }
    4 = @{ "B"=1.02; "C"=1.049453273958411; "D"=1.05343456307037; "E"=1.062167578651376; "F"=1.067811515222576; "I"=1.040619581662674; "J"=1.05392061727131; "K"=1.055875222631369; "L"=1.064587120909269; "M"=1.070217607183708; "N"=1.05541730543422 }
    5 = @{ "B"=1.02; "C"=1.04973388334357; "D"=1.053645146899111; "E"=1.062431271615628; "F"=1.068066733636945; "I"=1.040667821000523; "J"=1.05410513032445; "K"=1.056034471934684; "L"=1.064799823773029; "M"=1.070422154894644; "N"=1.055602080517059 }
    6 = @{ "B"=1.02; "C"=1.049781006596055; "D"=1.053680508806749; "E"=1.062475558365732; "F"=1.068109595351838; "I"=1.040675903756934; "J"=1.054136110087102; "K"=1.056061205294819; "L"=1.064835540326022; "M"=1.070456499941387; "N"=1.055633104274528 }
    7 = @{ "B"=1.02; "C"=1.04945702295825; "D"=1.053437376635902; "E"=1.062171101355601; "F"=1.067814924831516; "I"=1.040620227367195; "J"=1.053923082796217; "K"=1.055877350882814; "L"=1.064589962866682; "M"=1.070220340320301; "N"=1.055419774460455 }
    8 = @{ "B"=1.02; "C"=1.04810394300343; "D"=1.052421704432181; "E"=1.06090019838156; "F"=1.066584621397005; "I"=1.040385074959277; "J"=1.053032554349676; "K"=1.055108108433642; "L"=1.063563901544337; "M"=1.06923332309236; "N"=1.054527981361375 }
    9 = @{ "B"=1.02; "C"=1.045723899110219; "D"=1.050634225068721; "E"=1.058667098561658; "F"=1.064421913683467; "I"=1.039961560091639; "J"=1.051462968483635; "K"=1.053749791674845; "L"=1.061757464097659; "M"=1.067494480291669; "N"=1.05295616650337 }
    10 = @{ "B"=1.02; "C"=1.044140097660557; "D"=1.04944416283365; "E"=1.057182725712136; "F"=1.062983701098078; "I"=1.039673093646718; "J"=1.050416379700191; "K"=1.05284240480785; "L"=1.060554309786541; "M"=1.06633558164056; "N"=1.051908091443809 }
    11 = @{ "B"=1.02; "C"=1.043454987030966; "D"=1.048929244055632; "E"=1.056541022690256; "F"=1.062361804957284; "I"=1.039546737531743; "J"=1.04996315720868; "K"=1.052449069349509; "L"=1.060033611105493; "M"=1.065833855381651; "N"=1.051454225324373 }
    12 = @{ "B"=1.02; "C"=1.043200609914138; "D"=1.048738039383675; "E"=1.056302822864991; "F"=1.062130935214009; "I"=1.039499586062644; "J"=1.049794804467278; "K"=1.052302903066845; "L"=1.059840242569095; "M"=1.065647505477849; "N"=1.051285633502769 }
    13 = @{ "B"=1.02; "C"=1.043255169979677; "D"=1.048779050743415; "E"=1.056353910404692; "F"=1.062180451650021; "I"=1.039509710034331; "J"=1.049830916950104; "K"=1.05233425911357; "L"=1.05988171885669; "M"=1.065687477536785; "N"=1.051321797269462 }
    14 = @{ "B"=1.02; "C"=1.043433958025729; "D"=1.048913437806081; "E"=1.056521329804955; "F"=1.062342718530661; "I"=1.039542844404725; "J"=1.049949241231483; "K"=1.052436988499726; "L"=1.060017626337654; "M"=1.065818451358634; "N"=1.051440289584892 }
    15 = @{ "B"=1.02; "C"=1.043544128966275; "D"=1.048996245973292; "E"=1.056624503288999; "F"=1.062442713792227; "I"=1.03956323081717; "J"=1.050022143987704; "K"=1.052500275009989; "L"=1.060101369030163; "M"=1.065899150476992; "N"=1.05151329587139 }
    16 = @{ "B"=1.02; "C"=1.044185580683685; "D"=1.04947834452349; "E"=1.057225335427057; "F"=1.06302499249363; "I"=1.039681449007241; "J"=1.050446457743718; "K"=1.052868500172403; "L"=1.060588872707808; "M"=1.066368881414793; "N"=1.051938212201606 }
    17 = @{ "B"=1.02; "C"=1.044588130325335; "D"=1.049780856328345; "E"=1.057602500731901; "F"=1.063390471371668; "I"=1.039755216701104; "J"=1.050712607614276; "K"=1.053099363249481; "L"=1.060894745102395; "M"=1.066663554461454; "N"=1.052204740035509 }
    18 = @{ "B"=1.02; "C"=1.044822996865787; "D"=1.049957343518313; "E"=1.057822595173935; "F"=1.063603731866333; "I"=1.039798104387739; "J"=1.050867844170949; "K"=1.053233980143864; "L"=1.061073181781476; "M"=1.066835440358961; "N"=1.052360197045898 }
    19 = @{ "B"=1.02; "C"=1.044903091492812; "D"=1.050017527361058; "E"=1.057897658674267; "F"=1.063676462204244; "I"=1.039812704262513; "J"=1.05092077509891; "K"=1.05327987394114; "L"=1.06113402854572; "M"=1.066894050354604; "N"=1.052413203141848 }
    20 = @{ "B"=1.02; "C"=1.044544933713576; "D"=1.049748395821791; "E"=1.057562024082523; "F"=1.063351250369815; "I"=1.039747316578933; "J"=1.050684052677324; "K"=1.053074598126861; "L"=1.060861925143704; "M"=1.066631937988675; "N"=1.052176144547271 }
    21 = @{ "B"=1.02; "C"=1.043381306542389; "D"=1.048873862525708; "E"=1.05647202461577; "F"=1.062294931376777; "I"=1.039533093150562; "J"=1.049914397832442; "K"=1.052406739015969; "L"=1.059977603802178; "M"=1.065779882450973; "N"=1.051405396704229 }
    22 = @{ "B"=1.02; "C"=1.042650286852994; "D"=1.048324351144206; "E"=1.055787608963607; "F"=1.061631535368729; "I"=1.039397146029917; "J"=1.04943045247912; "K"=1.051986459696608; "L"=1.059421840075926; "M"=1.065244240250655; "N"=1.050920764092995 }
    23 = @{ "B"=1.02; "C"=1.043037757306455; "D"=1.048615624702834; "E"=1.056150344043047; "F"=1.061983142351601; "I"=1.039469333124131; "J"=1.049687003949938; "K"=1.052209292511233; "L"=1.059716437523929; "M"=1.065528186541913; "N"=1.051177679896329 }
    24 = @{ "B"=1.02; "C"=1.044564452198631; "D"=1.049763063212687; "E"=1.057580313428222; "F"=1.063368972396386; "I"=1.039750886735948; "J"=1.050696955437077; "K"=1.053085788548452; "L"=1.060876754988518; "M"=1.066646224086845; "N"=1.052189065630423 }
    25 = @{ "B"=1.02; "C"=1.046338689776567; "D"=1.051096055936068; "E"=1.059243644370941; "F"=1.064980397143573; "I"=1.040072130019935; "J"=1.051868782252273; "K"=1.054101277445917; "L"=1.06222427485796; "M"=1.067943958782371; "N"=1.053362556574107 }
}

foreach ($rowKey in $data.Keys) {
    $rowVals = $data[$rowKey]
    foreach ($colKey in $rowVals.Keys) {
        $ws.Range("$colKey$rowKey").Value = $rowVals[$colKey]
    }
}

Write-Output "Updated $($data.Count) rows of vm_pu results for 380 kV case"